$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$rng = $ws1.Range("A1:L21")
$rng.RemoveDuplicates(@(11,12), [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$sortRange = $ws1.Range("A2:L11")
$ws1.Sort.SortFields.Clear()
$ws1.Sort.SortFields.Add($ws1.Range("K2:K11"))
$ws1.Sort.SetRange($sortRange)
$ws1.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws1.Sort.Apply()
$ws1.Range("A12:L21").EntireRow.Delete()
